# Update dashboards - 2025-12-18
# Applies the scraped-data refresh described in the commit diff to the
# "Aguilar Prototype" worksheet.
#
# Notes on "blank" cells: throughout this sheet, data points that have
# no reading yet are stored as a numeric cell holding 0 (<v/>), not as
# a truly empty cell - e.g. S29/T29 in the untouched original. We keep
# that same convention when clearing Q/R/T/U cells below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 7 - GDP Nowcast present value
# ---------------------------------------------------------------------
$ws.Range("F7").Value = 0.1966919365954514

# ---------------------------------------------------------------------
# Row 13 - UI Initial Claims (ICSA): new latest date, values shift with
# a new "Present" reading pushing the rest of the lag columns along.
# ---------------------------------------------------------------------
$ws.Range("N13").Value = 45999
$ws.Range("Q13").Value = 224000
$ws.Range("R13").Value = 237000
$ws.Range("S13").Value = 192000
$ws.Range("T13").Value = 217000
$ws.Range("U13").Value = 222000

# ---------------------------------------------------------------------
# Row 14 - UI Continuing Claims (CCSA)
# ---------------------------------------------------------------------
$ws.Range("N14").Value = 45992
$ws.Range("Q14").Value = 1897000
$ws.Range("R14").Value = 1830000
$ws.Range("S14").Value = 1937000
$ws.Range("T14").Value = 1944000
$ws.Range("U14").Value = 1953000

# ---------------------------------------------------------------------
# Row 18 - CPI M/M % Delta: latest date moves forward a month, which
# also flips the "Latest Date" cell to the highlighted (yellow) style
# used for the most-recent release (style index 48, same as N13).
# ---------------------------------------------------------------------
$ws.Range("N13").Copy()
$ws.Range("N18").PasteSpecial(-4122)
$ws.Range("N18").Value = 45962
$ws.Range("Q18").Value = 0
$ws.Range("R18").Value = 0
$ws.Range("S18").Value = 0.00310486015759337
$ws.Range("T18").Value = 0.003824519141221616
$ws.Range("U18").Value = 0.00196578538102643

# ---------------------------------------------------------------------
# Row 19 - CPI Y/Y % Delta
# ---------------------------------------------------------------------
$ws.Range("N13").Copy()
$ws.Range("N19").PasteSpecial(-4122)
$ws.Range("N19").Value = 45962
$ws.Range("Q19").Value = 0.0271196938527219
$ws.Range("R19").Value = 0
$ws.Range("S19").Value = 0.03022699626172379
$ws.Range("T19").Value = 0.02939219624933549
$ws.Range("U19").Value = 0.02731801279475463

# ---------------------------------------------------------------------
# Row 20 - Core CPI M/M % Delta
# ---------------------------------------------------------------------
$ws.Range("N13").Copy()
$ws.Range("N20").PasteSpecial(-4122)
$ws.Range("N20").Value = 45962
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("S20").Value = 0.002271121582325675
$ws.Range("T20").Value = 0.003459544325982167
$ws.Range("U20").Value = 0.003223443223443256

# ---------------------------------------------------------------------
# Row 21 - Core CPI Y/Y % Delta
# ---------------------------------------------------------------------
$ws.Range("N13").Copy()
$ws.Range("N21").PasteSpecial(-4122)
$ws.Range("N21").Value = 45962
$ws.Range("Q21").Value = 0.02618878615332623
$ws.Range("R21").Value = 0
$ws.Range("S21").Value = 0.03025542724453378
$ws.Range("T21").Value = 0.03112190821006822
$ws.Range("U21").Value = 0.03048602684576389

# ---------------------------------------------------------------------
# Row 29 - 5yr, 5yr Forward inflation expectation
# ---------------------------------------------------------------------
$ws.Range("N29").Value = 46008
$ws.Range("Q29").Value = 2.22
$ws.Range("R29").Value = 2.21
$ws.Range("S29").Value = 2.21
$ws.Range("U29").Value = 0

# ---------------------------------------------------------------------
# Row 30 - 10yr TIPS breakeven
# ---------------------------------------------------------------------
$ws.Range("N30").Value = 46008
$ws.Range("Q30").Value = 2.24
$ws.Range("R30").Value = 2.23
$ws.Range("S30").Value = 2.25
$ws.Range("U30").Value = 0

# ---------------------------------------------------------------------
# Row 46 - Exports header row: "Latest Date" style reverts to the
# non-highlighted style (47), value unchanged.
# ---------------------------------------------------------------------
$ws.Range("N22").Copy()
$ws.Range("C46").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 47 - Exports M/M % Delta SA; FFR lag data
# ---------------------------------------------------------------------
$ws.Range("N22").Copy()
$ws.Range("C47").PasteSpecial(-4122)
$ws.Range("N47").Value = 46007

# ---------------------------------------------------------------------
# Row 48 - Imports; 2y UST
# ---------------------------------------------------------------------
$ws.Range("N22").Copy()
$ws.Range("C48").PasteSpecial(-4122)
$ws.Range("N48").Value = 46007
$ws.Range("Q48").Value = 3.48
$ws.Range("R48").Value = 3.51
$ws.Range("T48").Value = 0

# ---------------------------------------------------------------------
# Row 49 - Imports M/M % Delta SA; 5y UST
# ---------------------------------------------------------------------
$ws.Range("N22").Copy()
$ws.Range("C49").PasteSpecial(-4122)
$ws.Range("N49").Value = 46007
$ws.Range("Q49").Value = 3.69
$ws.Range("R49").Value = 3.73
$ws.Range("T49").Value = 0
$ws.Range("U49").Value = 3.75

# ---------------------------------------------------------------------
# Row 50 - Trade Balance; 10y UST
# ---------------------------------------------------------------------
$ws.Range("N22").Copy()
$ws.Range("C50").PasteSpecial(-4122)
$ws.Range("N50").Value = 46007
$ws.Range("Q50").Value = 4.15
$ws.Range("R50").Value = 4.18
$ws.Range("T50").Value = 0
$ws.Range("U50").Value = 4.19

# ---------------------------------------------------------------------
# Row 51 - Trade Balance M/M % Delta SA; 30y Mortgage
# ---------------------------------------------------------------------
$ws.Range("N22").Copy()
$ws.Range("C51").PasteSpecial(-4122)
$ws.Range("N22").Copy()
$ws.Range("N51").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Row 52 - BAA corporate bond yield
# ---------------------------------------------------------------------
$ws.Range("N52").Value = 46007
$ws.Range("Q52").Value = 5.93
$ws.Range("R52").Value = 5.95
$ws.Range("T52").Value = 0
$ws.Range("U52").Value = 5.95

$excel.CutCopyMode = $false
